# Update the dSF (column F) values for the listed rows, per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = -8
$ws.Range("F12").Value = -5
$ws.Range("F15").Value = -3
